# Group 4 - Final Presentation (Term Project).pptx
# Slide 1 title: "CSE 4288 – Term Project" <line break> "Progress Presentation"
#              -> "CSE 4288 – Term Project" <line break> "Final Presentation"
#
# The existing line break (a:br) must stay exactly where it is, so we only
# touch the two words after it: "Progress" becomes "Final " and the leading
# space on " Presentation" is dropped (since "Final " now supplies it).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# Replace "Progress" with "Final "
$text = $tr.Text
$progressStart = $text.IndexOf("Progress")
$wordRange = $tr.Characters($progressStart + 1, 8)
$wordRange.Text = "Final "

# Replace " Presentation" (leading space) with "Presentation"
$text = $tr.Text
$presStart = $text.IndexOf(" Presentation")
$presRange = $tr.Characters($presStart + 1, 13)
$presRange.Text = "Presentation"
